# Swap the roster rows for "Malcolm Brogdon" (row 6) and "Jaylen Brown" (row 7).
# Column A (the 0-based "No." index) stays put; columns B..K (No, Player, Pos, Ht,
# Wt, Birth Date, Unnamed: 6, Exp, College, bbref url) swap between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 6
$row2 = 7
$firstCol = 2   # column B
$lastCol = 11   # column K

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)

    $v1 = $cell1.Value()
    $v2 = $cell2.Value()

    $cell1.Value = $v2
    $cell2.Value = $v1
}
